# "removed previous demo data and added one record for each model"
#
# The demo sheet previously held 4 sample rows (rows 2-5). We wipe that out
# and leave a single fresh record in row 2, matching the pattern used for
# the other demo-data sheets (one row of sample data per model).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122
# Row 5 (A5) already carries the sheet's "plain" default style (no explicit
# number format / wrap). Stamp that same default style onto B5:C5 and onto
# A2 before we touch any values, so that once the values are cleared below,
# every one of those cells collapses back to the implicit default style
# (i.e. the cell node itself disappears rather than lingering as an empty,
# explicitly-styled cell).
$ws.Range("A5").Copy()
$ws.Range("B5:C5").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Drop all the old demo rows (2-5).
$ws.Range("A2:C5").ClearContents()

# ... and add back exactly one new record.
$ws.Range("A2").Value = "200 mg"
$ws.Range("B2").Value = "Tablet"
$ws.Range("C2").Value = 44227

# Leave the selection where the user's cursor ended up after typing the row.
$ws.Range("A6").Select() | Out-Null
